# "Nädal 2" ("Week 2") is already the active sheet/tab in this workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 of the log (worksheet row 17) was filled in with a new entry:
#   Date  = 08/02/2020
#   Start = 19:30 (0.8125 as an Excel time fraction)
#   Activity = "kood"
#   Comments = "RP with EFCore, juhendi järgi"
$ws.Range("B17").Value = (Get-Date -Year 2020 -Month 2 -Day 8).Date
$ws.Range("C17").Value = 0.8125
$ws.Range("G17").Value = "kood"
$ws.Range("H17").Value = "RP with EFCore, juhendi järgi"

# Column C (Start time) got wider to fit the newly entered values.
$ws.Columns("C").ColumnWidth = 15

# The active selection moved to H18.
$ws.Range("H18").Select() | Out-Null
